$d = $word.ActiveDocument

# 1) Intro paragraph: clarify startup state for prepackaged VM
$d.Content.Find.Execute(
    "Boot your Linux system or VM, log in, and then open a terminal window and start the lab:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Boot your Linux system or VM.  If necessary, log in and then open a terminal window and cd to the labtainer/labtainer-student directory.  The pre-packaged Labtainer VM will start with such a terminal open for you.   Then start the lab:",
    2) | Out-Null

# 2) Remove the now-redundant "cd labtainer/labtainer-student" command line
$d.Content.Find.Execute(
    "cd labtainer/labtainer-student",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2) | Out-Null

# 3) Rename start.py -> labtainer command
$d.Content.Find.Execute(
    "start.py sys-log",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "labtainer sys-log",
    2) | Out-Null

# 4) Drop the "stretch the terminal" sentence before "Note the terminal displays..."
$d.Content.Find.Execute(
    "It may help to stretch the resulting bash terminal window to the right to provide for more output space.  Note the terminal displays the paths to two files on your Linux host: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Note the terminal displays the paths to two files on your Linux host: ",
    2) | Out-Null

# 5) stop.py -> stoplab (first, quoted, bold context)
$d.Content.Find.Execute(
    [char]8220 + "stop.py" + [char]8221 + " to stop the lab for the last time.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]8220 + "stoplab" + [char]8221 + " to stop the lab for the last time.",
    2) | Out-Null

# 6) Merge split runs "/etc/" + "rsyslog.d/50-default" + ".conf" (first occurrence)
$d.Content.Find.Execute(
    "/etc/rsyslog.d/50-default.conf",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "/etc/rsyslog.d/50-default.conf",
    2) | Out-Null

# 7) stop.py sys-log -> stoplab sys-log
$d.Content.Find.Execute(
    "stop.py sys-log",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "stoplab sys-log",
    2) | Out-Null

# 8) "./stop.py" -> "stoplab"
$d.Content.Find.Execute(
    [char]8220 + "./stop.py" + [char]8221 + ".",
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]8220 + "stoplab" + [char]8221 + ".",
    2) | Out-Null

# 9) Appendix table: bump left cell margin from 143dxa (7.15pt) to 153dxa (7.65pt)
$d.Tables.Item(1).LeftPadding = 7.65
